$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4499.5
$ws.Range("I40").Value = 5133.3335
$ws.Range("J40").Value = 2598
$ws.Range("K40").Value = 5133.3335
$ws.Range("L40").Value = 2598
$ws.Range("M40").Value = -4958.3335
$ws.Range("N40").Value = -2948
$ws.Range("H44").Value = 4406.8184
$ws.Range("I44").Value = 4406.8184
$ws.Range("K44").Value = 4406.8184
$ws.Range("M44").Value = -3944.8184
$ws.Range("H62").Value = 5570.5625
$ws.Range("I62").Value = 4023
$ws.Range("K62").Value = 4023
$ws.Range("M62").Value = -3399
$ws.Range("H65").Value = 5570.5625
$ws.Range("I65").Value = 4023
$ws.Range("K65").Value = 20115
$ws.Range("M65").Value = -16995
$ws.Range("H132").Value = 2549.4773
$ws.Range("I132").Value = 2207.8918
$ws.Range("J132").Value = 4355
$ws.Range("K132").Value = 6623.6754
$ws.Range("L132").Value = 13065
$ws.Range("M132").Value = -4093.6754
$ws.Range("N132").Value = -18125
$ws.Range("H135").Value = 50001956
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 50001956
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 450017604
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -450022674
$ws.Range("H137").Value = 3467.3225
$ws.Range("I137").Value = 3563.9443
$ws.Range("J137").Value = 3333.5386
$ws.Range("K137").Value = 10691.8329
$ws.Range("L137").Value = 10000.6158
$ws.Range("M137").Value = -8141.832900000001
$ws.Range("N137").Value = -15100.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16672957
$ws.Range("I32").Value = 20411548
$ws.Range("J32").Value = 19236.092
$ws.Range("K32").Value = 20411548
$ws.Range("L32").Value = 19236.092
$ws.Range("M32").Value = -20411261
$ws.Range("N32").Value = -19810.092
$ws.Range("H61").Value = 83339210
$ws.Range("I61").Value = 142860720
$ws.Range("J61").Value = 9101.6
$ws.Range("K61").Value = 142860720
$ws.Range("L61").Value = 9101.6
$ws.Range("M61").Value = -142860508
$ws.Range("N61").Value = -9525.6
$ws.Range("H69").Value = 200459
$ws.Range("J69").Value = 200459
$ws.Range("L69").Value = 200459
$ws.Range("N69").Value = -201957
$ws.Range("H72").Value = 200459
$ws.Range("J72").Value = 200459
$ws.Range("L72").Value = 601377
$ws.Range("N72").Value = -608865
$ws.Range("H74").Value = 47621348
$ws.Range("I74").Value = 58825396
$ws.Range("J74").Value = 4128
$ws.Range("K74").Value = 58825396
$ws.Range("L74").Value = 4128
$ws.Range("M74").Value = -58824522
$ws.Range("N74").Value = -5876
$ws.Range("H77").Value = 47621348
$ws.Range("I77").Value = 58825396
$ws.Range("J77").Value = 4128
$ws.Range("K77").Value = 294126980
$ws.Range("L77").Value = 20640
$ws.Range("M77").Value = -294122612
$ws.Range("N77").Value = -29376
$ws.Range("H122").Value = 2947.111
$ws.Range("I122").Value = 1895.2106
$ws.Range("K122").Value = 5685.6318
$ws.Range("M122").Value = -3235.6318
$ws.Range("H136").Value = 83339210
$ws.Range("I136").Value = 142860720
$ws.Range("J136").Value = 9101.6
$ws.Range("K136").Value = 428582160
$ws.Range("L136").Value = 27304.8
$ws.Range("M136").Value = -428579610
$ws.Range("N136").Value = -32404.8
$ws.Range("H139").Value = 78571.664
$ws.Range("J139").Value = 78571.664
$ws.Range("L139").Value = 78571.664
$ws.Range("N139").Value = -88851.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 50000
$ws.Range("J76").Value = 50000
$ws.Range("L76").Value = 50000
$ws.Range("N76").Value = -50630
$ws.Range("H79").Value = 50000
$ws.Range("J79").Value = 50000
$ws.Range("L79").Value = 50000
$ws.Range("N79").Value = -52184
$ws.Range("H94").Value = 1091.9656
$ws.Range("I94").Value = 835.6957
$ws.Range("J94").Value = 2074.3333
$ws.Range("K94").Value = 835.6957
$ws.Range("L94").Value = 2074.3333
$ws.Range("M94").Value = -384.6957
$ws.Range("N94").Value = -2976.3333
$ws.Range("H105").Value = 11706.3
$ws.Range("I105").Value = 15766.286
$ws.Range("K105").Value = 15766.286
$ws.Range("M105").Value = -14019.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23814064
$ws.Range("I31").Value = 3551.963
$ws.Range("J31").Value = 66672984
$ws.Range("K31").Value = 3551.963
$ws.Range("L31").Value = 66672984
$ws.Range("M31").Value = -3256.963
$ws.Range("N31").Value = -66673574
$ws.Range("H34").Value = 23814064
$ws.Range("I34").Value = 3551.963
$ws.Range("J34").Value = 66672984
$ws.Range("K34").Value = 3551.963
$ws.Range("L34").Value = 66672984
$ws.Range("M34").Value = -3349.963
$ws.Range("N34").Value = -66673388
$ws.Range("H58").Value = 2795.4062
$ws.Range("I58").Value = 2534.7778
$ws.Range("J58").Value = 4202.8
$ws.Range("K58").Value = 2534.7778
$ws.Range("L58").Value = 4202.8
$ws.Range("M58").Value = -2331.7778
$ws.Range("N58").Value = -4608.8
$ws.Range("H107").Value = 2036.75
$ws.Range("I107").Value = 1129.8
$ws.Range("K107").Value = 1129.8
$ws.Range("M107").Value = 790.2
$ws.Range("H122").Value = 1506.5294
$ws.Range("I122").Value = 1339.8148
$ws.Range("K122").Value = 4019.4444
$ws.Range("M122").Value = -1569.4444
$ws.Range("H132").Value = 2919.6428
$ws.Range("I132").Value = 2122.261
$ws.Range("J132").Value = 6587.6
$ws.Range("K132").Value = 6366.782999999999
$ws.Range("L132").Value = 19762.8
$ws.Range("M132").Value = -3836.782999999999
$ws.Range("N132").Value = -24822.8
$ws.Range("H134").Value = 1288.1111
$ws.Range("I134").Value = 1185.75
$ws.Range("J134").Value = 2107
$ws.Range("K134").Value = 3557.25
$ws.Range("L134").Value = 6321
$ws.Range("M134").Value = -1022.25
$ws.Range("N134").Value = -11391
$ws.Range("H136").Value = 2795.4062
$ws.Range("I136").Value = 2534.7778
$ws.Range("J136").Value = 4202.8
$ws.Range("K136").Value = 7604.3334
$ws.Range("L136").Value = 12608.4
$ws.Range("M136").Value = -5054.3334
$ws.Range("N136").Value = -17708.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 5709.9165
$ws.Range("I75").Value = 7544.7144
$ws.Range("J75").Value = 3141.2
$ws.Range("K75").Value = 22634.1432
$ws.Range("L75").Value = 9423.599999999999
$ws.Range("M75").Value = -21636.1432
$ws.Range("N75").Value = -11419.6
$ws.Range("H78").Value = 5709.9165
$ws.Range("I78").Value = 7544.7144
$ws.Range("J78").Value = 3141.2
$ws.Range("K78").Value = 67902.4296
$ws.Range("L78").Value = 28270.8
$ws.Range("M78").Value = -62910.4296
$ws.Range("N78").Value = -38254.8
$ws.Range("H93").Value = 750
$ws.Range("J93").Value = 750
$ws.Range("L93").Value = 2250
$ws.Range("N93").Value = -5994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 27280218
$ws.Range("I126").Value = 25010854
$ws.Range("K126").Value = 75032562
$ws.Range("M126").Value = -75030092
$ws.Range("H132").Value = 3461.08
$ws.Range("I132").Value = 3207.7058
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 9623.117400000001
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -7093.117400000001
$ws.Range("N132").Value = -17058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3859.4443
$ws.Range("I40").Value = 4091.875
$ws.Range("K40").Value = 4091.875
$ws.Range("M40").Value = -3955.875
$ws.Range("H46").Value = 1154.6305
$ws.Range("I46").Value = 663.74286
$ws.Range("J46").Value = 2716.5454
$ws.Range("K46").Value = 663.74286
$ws.Range("L46").Value = 2716.5454
$ws.Range("M46").Value = -475.74286
$ws.Range("N46").Value = -3092.5454
$ws.Range("H136").Value = 2072.131
$ws.Range("I136").Value = 1532.7273
$ws.Range("J136").Value = 7016.6665
$ws.Range("K136").Value = 4598.1819
$ws.Range("L136").Value = 21049.9995
$ws.Range("M136").Value = -2048.1819
$ws.Range("N136").Value = -26149.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3355.0303
$ws.Range("I126").Value = 4260.875
$ws.Range("K126").Value = 12782.625
$ws.Range("M126").Value = -10312.625
